# Generate Report for Handback
#
# A second file (7e8d9ea9-e727-4f68-90c9-4c120f6ef1a0.md) has now been
# handed back alongside the one already on the report. Re-running the
# report also refreshed the identifying hash/UUID used for the file that
# was already present (285d2b4f-... -> 4fe2e3d8-...) together with its
# handoff/handback timestamps. This adds a row to each of the three
# sheets (Overview, zh-cn, de-de) and updates the existing row's
# volatile identifiers/timestamps.

$wb = $excel.ActiveWorkbook
$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsO = $wb.Worksheets.Item("Overview")

# Row 2 - refresh UUID + hyperlink + timestamp of the already-handed-back file.
$wsO.Range("A2").Value = "4fe2e3d8-84c7-440c-baa1-3649045a245a.md"
$wsO.Range("B2").Hyperlinks.Delete()
$wsO.Range("B2").Value = "e2e\4fe2e3d8-84c7-440c-baa1-3649045a245a.md"
$wsO.Hyperlinks.Add($wsO.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/058f1c5ce4715f0482d3e340d18560d644410305/e2e/4fe2e3d8-84c7-440c-baa1-3649045a245a.md", "", "", "e2e\4fe2e3d8-84c7-440c-baa1-3649045a245a.md")
$wsO.Range("C2").Value = ".md"
$wsO.Range("E2").Value = "Handed back: in sync with en-US"
$wsO.Range("F2").Value = "Handed back: in sync with en-US"
$wsO.Range("G2").Value = "2016-09-04 17:06:57"
$wsO.Range("G2").NumberFormat = $dateFmt

# Row 3 - newly handed-back file.
$wsO.Range("A3").Value = "7e8d9ea9-e727-4f68-90c9-4c120f6ef1a0.md"
$wsO.Range("B3").Value = "e2e\7e8d9ea9-e727-4f68-90c9-4c120f6ef1a0.md"
$wsO.Hyperlinks.Add($wsO.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/058f1c5ce4715f0482d3e340d18560d644410305/e2e/7e8d9ea9-e727-4f68-90c9-4c120f6ef1a0.md", "", "", "e2e\7e8d9ea9-e727-4f68-90c9-4c120f6ef1a0.md")
$wsO.Range("C3").Value = ".md"
$wsO.Range("E3").Value = "Handed back: in sync with en-US"
$wsO.Range("F3").Value = "Handed back: in sync with en-US"
$wsO.Range("G3").Value = "2016-09-04 17:06:57"
$wsO.Range("G3").NumberFormat = $dateFmt

$loO = $wsO.ListObjects.Item(1)
$loO.Resize($wsO.Range("A1:G3"))

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZ = $wb.Worksheets.Item("zh-cn")

# Row 2 - refresh UUID/hash/hyperlinks/timestamps of the existing file.
$wsZ.Range("A2").Hyperlinks.Delete()
$wsZ.Range("A2").Value = "4fe2e3d8-84c7-440c-baa1-3649045a245a.md"
$wsZ.Hyperlinks.Add($wsZ.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/058f1c5ce4715f0482d3e340d18560d644410305/e2e/4fe2e3d8-84c7-440c-baa1-3649045a245a.md", "", "", "4fe2e3d8-84c7-440c-baa1-3649045a245a.md")
$wsZ.Range("B2").Value = ".md"
$wsZ.Range("C2").Value = "Handed back: in sync with en-US"
$wsZ.Range("D2").Value = "e2e"
$wsZ.Range("E2").Value = "ht"
$wsZ.Range("F2").Value = "False"
$wsZ.Range("G2").Value = "4fe2e3d8-84c7-440c-baa1-3649045a245a.d940987b3ddca58a451eca05e8623f7d668d62d2.zh-cn.xlf"
$wsZ.Range("H2").Value = "2016-09-04 17:06:52"
$wsZ.Range("H2").NumberFormat = $dateFmt
$wsZ.Range("I2").Hyperlinks.Delete()
$wsZ.Range("I2").Value = "4fe2e3d8-84c7-440c-baa1-3649045a245a.md"
$wsZ.Hyperlinks.Add($wsZ.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/0ab278192a184ba5a5a3c081119b615c5be7ce2e/e2e/4fe2e3d8-84c7-440c-baa1-3649045a245a.md", "", "", "4fe2e3d8-84c7-440c-baa1-3649045a245a.md")
$wsZ.Range("J2").Value = "4fe2e3d8-84c7-440c-baa1-3649045a245a.d940987b3ddca58a451eca05e8623f7d668d62d2.zh-cn.xlf"
$wsZ.Range("K2").Value = "2016-09-04 17:07:15"
$wsZ.Range("K2").NumberFormat = $dateFmt
$wsZ.Range("L2").Value = ""
$wsZ.Range("M2").Value = "True"
$wsZ.Range("N2").Value = ""
$wsZ.Range("O2").Value = "False"
$wsZ.Range("P2").Value = ""

# Row 3 - newly handed-back file.
$wsZ.Range("A3").Value = "7e8d9ea9-e727-4f68-90c9-4c120f6ef1a0.md"
$wsZ.Hyperlinks.Add($wsZ.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/058f1c5ce4715f0482d3e340d18560d644410305/e2e/7e8d9ea9-e727-4f68-90c9-4c120f6ef1a0.md", "", "", "7e8d9ea9-e727-4f68-90c9-4c120f6ef1a0.md")
$wsZ.Range("B3").Value = ".md"
$wsZ.Range("C3").Value = "Handed back: in sync with en-US"
$wsZ.Range("D3").Value = "e2e"
$wsZ.Range("E3").Value = "ht"
$wsZ.Range("F3").Value = "True"
$wsZ.Range("G3").Value = "7e8d9ea9-e727-4f68-90c9-4c120f6ef1a0.734709e17f10cda3c3eea1de08e49da228698c04.zh-cn.xlf"
$wsZ.Range("H3").Value = "2016-09-04 17:06:52"
$wsZ.Range("H3").NumberFormat = $dateFmt
$wsZ.Range("I3").Value = "7e8d9ea9-e727-4f68-90c9-4c120f6ef1a0.md"
$wsZ.Hyperlinks.Add($wsZ.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/0ab278192a184ba5a5a3c081119b615c5be7ce2e/e2e/7e8d9ea9-e727-4f68-90c9-4c120f6ef1a0.md", "", "", "7e8d9ea9-e727-4f68-90c9-4c120f6ef1a0.md")
$wsZ.Range("J3").Value = "7e8d9ea9-e727-4f68-90c9-4c120f6ef1a0.734709e17f10cda3c3eea1de08e49da228698c04.zh-cn.xlf"
$wsZ.Range("K3").Value = "2016-09-04 17:07:15"
$wsZ.Range("K3").NumberFormat = $dateFmt
$wsZ.Range("L3").Value = ""
$wsZ.Range("M3").Value = "True"
$wsZ.Range("N3").Value = ""
$wsZ.Range("O3").Value = "False"
$wsZ.Range("P3").Value = ""

$loZ = $wsZ.ListObjects.Item(1)
$loZ.Resize($wsZ.Range("A1:P3"))

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsD = $wb.Worksheets.Item("de-de")

# Row 2 - refresh UUID/hash/hyperlinks/timestamps of the existing file.
$wsD.Range("A2").Hyperlinks.Delete()
$wsD.Range("A2").Value = "4fe2e3d8-84c7-440c-baa1-3649045a245a.md"
$wsD.Hyperlinks.Add($wsD.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/058f1c5ce4715f0482d3e340d18560d644410305/e2e/4fe2e3d8-84c7-440c-baa1-3649045a245a.md", "", "", "4fe2e3d8-84c7-440c-baa1-3649045a245a.md")
$wsD.Range("B2").Value = ".md"
$wsD.Range("C2").Value = "Handed back: in sync with en-US"
$wsD.Range("D2").Value = "e2e"
$wsD.Range("E2").Value = "ht"
$wsD.Range("F2").Value = "False"
$wsD.Range("G2").Value = "4fe2e3d8-84c7-440c-baa1-3649045a245a.d940987b3ddca58a451eca05e8623f7d668d62d2.de-de.xlf"
$wsD.Range("H2").Value = "2016-09-04 17:06:57"
$wsD.Range("H2").NumberFormat = $dateFmt
$wsD.Range("I2").Hyperlinks.Delete()
$wsD.Range("I2").Value = "4fe2e3d8-84c7-440c-baa1-3649045a245a.md"
$wsD.Hyperlinks.Add($wsD.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/a45e98aeb6a4ac6a3f5ef9e7b2f42d03e64163f5/e2e/4fe2e3d8-84c7-440c-baa1-3649045a245a.md", "", "", "4fe2e3d8-84c7-440c-baa1-3649045a245a.md")
$wsD.Range("J2").Value = "4fe2e3d8-84c7-440c-baa1-3649045a245a.d940987b3ddca58a451eca05e8623f7d668d62d2.de-de.xlf"
$wsD.Range("K2").Value = "2016-09-04 17:07:23"
$wsD.Range("K2").NumberFormat = $dateFmt
$wsD.Range("L2").Value = ""
$wsD.Range("M2").Value = "True"
$wsD.Range("N2").Value = ""
$wsD.Range("O2").Value = "False"
$wsD.Range("P2").Value = ""

# Row 3 - newly handed-back file.
$wsD.Range("A3").Value = "7e8d9ea9-e727-4f68-90c9-4c120f6ef1a0.md"
$wsD.Hyperlinks.Add($wsD.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/058f1c5ce4715f0482d3e340d18560d644410305/e2e/7e8d9ea9-e727-4f68-90c9-4c120f6ef1a0.md", "", "", "7e8d9ea9-e727-4f68-90c9-4c120f6ef1a0.md")
$wsD.Range("B3").Value = ".md"
$wsD.Range("C3").Value = "Handed back: in sync with en-US"
$wsD.Range("D3").Value = "e2e"
$wsD.Range("E3").Value = "ht"
$wsD.Range("F3").Value = "True"
$wsD.Range("G3").Value = "7e8d9ea9-e727-4f68-90c9-4c120f6ef1a0.734709e17f10cda3c3eea1de08e49da228698c04.de-de.xlf"
$wsD.Range("H3").Value = "2016-09-04 17:06:57"
$wsD.Range("H3").NumberFormat = $dateFmt
$wsD.Range("I3").Value = "7e8d9ea9-e727-4f68-90c9-4c120f6ef1a0.md"
$wsD.Hyperlinks.Add($wsD.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/a45e98aeb6a4ac6a3f5ef9e7b2f42d03e64163f5/e2e/7e8d9ea9-e727-4f68-90c9-4c120f6ef1a0.md", "", "", "7e8d9ea9-e727-4f68-90c9-4c120f6ef1a0.md")
$wsD.Range("J3").Value = "7e8d9ea9-e727-4f68-90c9-4c120f6ef1a0.734709e17f10cda3c3eea1de08e49da228698c04.de-de.xlf"
$wsD.Range("K3").Value = "2016-09-04 17:07:23"
$wsD.Range("K3").NumberFormat = $dateFmt
$wsD.Range("L3").Value = ""
$wsD.Range("M3").Value = "True"
$wsD.Range("N3").Value = ""
$wsD.Range("O3").Value = "False"
$wsD.Range("P3").Value = ""

$loD = $wsD.ListObjects.Item(1)
$loD.Resize($wsD.Range("A1:P3"))

Write-Output "handback-status.xlsx updated: added 7e8d9ea9-e727-4f68-90c9-4c120f6ef1a0.md"
